$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Stimulus": update the existing "Basic transfers" / "Write Strobe"
# rows and append two new rows for IDLE/BUSY transfer test cases, plus a
# trailing "Cross Feature" row.
# ---------------------------------------------------------------------------
$stim = $wb.Worksheets.Item("Stimulus")

$stim.Range("C5").Value = "Write Strobe"
$stim.Range("D5").Value = "Full word write: Test a 32-bit data bus with all strobes active (HWSTRB = 1111)."

$stim.Range("D6").Value = "Sparse Write: Test a 32-bit data bus with only selected strobes active (HWSTRB = 1010 updates only bytes 0 and 2)."

$stim.Range("D7").Value = "No Write (All Strobes Inactive)"

$stim.Range("D8").Value = "Cross feature with hsize (HSIZE = HALF_WORD, HWSTRB = 0001)"

$stim.Range("C9").Value = "IDLE Transfer"
$stim.Range("D9").Value = "Send an IDLE transfer and verify that no read/write occurs and subordinate responds with OKAY."

# New row 10: BUSY Transfer (clone formatting from row 9, then overwrite values)
$stim.Range("B9:D9").Copy($stim.Range("B10:D10"))
$stim.Range("B10").Value = 7
$stim.Range("C10").Value = "BUSY Transfer"
$stim.Range("D10").Value = "Send an BUSY transfer and verify that no read/write occurs and subordinate responds with OKAY."

# New row 11: Cross Feature (clone formatting from row 10, then overwrite values)
$stim.Range("B10:D10").Copy($stim.Range("B11:D11"))
$stim.Range("B11").Value = 8
$stim.Range("C11").Value = "Cross Feature"
$stim.Range("D11").Value = "Multiple Write + Multiple Read txns with hsize, haddr same for each set of txns, with other fields randomize."

# ---------------------------------------------------------------------------
# Sheet "Coverage": the coverage-plan tables (state_cp / state_write_cross)
# are emptied out, leaving only the blank formatted rows behind.
# ---------------------------------------------------------------------------
$cov = $wb.Worksheets.Item("Coverage")
$cov.Range("B4:D14").ClearContents()
$cov.Range("D22").Select()

# Leave "Stimulus" as the active sheet/tab with D10 selected, matching the
# saved view state.
$stim.Range("D10").Select()
